$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 81, column A (date) value
$ws.Range("A81").Value = 45460.2916666667

# Add new row 82 with data
$ws.Range("A82").Value = 45461.6023263889
$ws.Range("B82").Value = 3000
$ws.Range("C82").Value = 2.95000004768372
$ws.Range("D82").Value = 2.95000004768372
$ws.Range("E82").Value = 2.95000004768372
$ws.Range("F82").Value = 2.95000004768372
$ws.Range("G82").NumberFormat = "@"
$ws.Range("G82").Value = "2.95000004768372"
$ws.Range("G82").ClearFormats()
$ws.Range("H82").Value = "ESPE.MI"

# Copy the style (date format) from A81 to A82
$ws.Range("A81").Copy()
$ws.Range("A82").PasteSpecial(-4122)  # xlPasteFormats
